# Generate Report for Handoff
# Updates the localization-status workbook to reflect that "b.md" has
# been handed off again (new handoff xliff files produced), while its
# previous handback turned out to be stale.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6353f3f93a63b16e4b7a1dd7ef7223a8784e4488/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/20b02bcba45a64e16ac11d1c26fe5d2e0da72c65/e2e/b.md."

# ---------------------------------------------------------------
# "Overview" sheet: b.md row (row 3)
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-25 08:40:20"

# ---------------------------------------------------------------
# "zh-cn" sheet: b.md row (row 3)
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
# Leading apostrophe forces text (rather than boolean) entry, as in
# the Excel UI; reset the style afterwards so no quote-prefix / extra
# style is left behind on the cell.
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("F3").Style = "Normal"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-25 08:40:07"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# ---------------------------------------------------------------
# "de-de" sheet: b.md row (row 3)
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("F3").Style = "Normal"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-25 08:40:20"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
